$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 9
$ws.Range("F3").Value = -7
$ws.Range("F4").Value = -7
$ws.Range("F5").Value = -5
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = -2
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 10
$ws.Range("F13").Value = 7
$ws.Range("F15").Value = -1
$ws.Range("F17").Value = -2
$ws.Range("F20").Value = -6
$ws.Range("F21").Value = -6
$ws.Range("F23").Value = -7
$ws.Range("F26").Value = -2
$ws.Range("F28").Value = -7
$ws.Range("F39").Value = -3
$ws.Range("F50").Value = -3
